{"js": "// Builds a minimal Flat-OPC wrapper so we can hand exact OOXML (with\n// precise <w:r>/<w:proofErr>/<w:bookmarkStart|End> markup) to\n// Range.insertOoxml / Paragraph.getRange().insertOoxml.\nfunction flatOpc(bodyXml) {\n  return '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n    '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst body = context.document.body;\n\n// --- Change 1 --------------------------------------------------------\n// \"Initialize \u2013 al[_GoBack]ways init before using. That means init\n// before starting threads that use it.\" loses the mid-word bookmark\n// split and becomes one clean run: \"Initialize \u2013 always init...\".\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet initPara = null;\nfor (const p of paras.items) {\n  if (p.text.indexOf(\"Initialize\") !== -1 && p.text.indexOf(\"before using\") !== -1) {\n    initPara = p;\n    break;\n  }\n}\n\nif (initPara) {\n  const fixedParaXml =\n    '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"39\"/></w:numPr></w:pPr>' +\n    '<w:r><w:t>Initialize</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> \\u2013 always </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> before using. That means </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> before starting threads that use it. </w:t></w:r>' +\n    '</w:p>';\n  initPara.getRange().insertOoxml(flatOpc(fixedParaXml), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- Change 2 --------------------------------------------------------\n// Add a new indented paragraph right after the WICED_WAIT_FOREVER line\n// explaining that the \"forever\" timeout is really ~50 days; the\n// _GoBack bookmark (freed up by change 1) now wraps the point right\n// before \"days)\".\nconst results = body.search(\"use WICED_WAIT_FOREVER if you don\\u2019t want a timeout\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const newParaXml =\n    '<w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr>' +\n    '<w:r><w:tab/><w:t xml:space=\"preserve\">(it\\u2019s really 2^32 </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>ms</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> which is just under 50 </w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r><w:t>days)</w:t></w:r>' +\n    '</w:p>';\n  results.items[0].insertOoxml(flatOpc(newParaXml), Word.InsertLocation.after);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Minimal Flat-OPC wrapper so we can hand Word exact OOXML (precise\n# <w:r>/<w:proofErr>/<w:bookmarkStart|End> markup) via Range.InsertXML.\nfunction New-FlatOpc([string]$bodyXml) {\n    return '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\">' +\n        '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n        '</Relationships></pkg:xmlData></pkg:part>' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyXml + '<w:sectPr/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$enDash = [char]8211\n$rsquo = [char]8217\n\n# --- Change 1 ----------------------------------------------------------\n# \"Initialize - al[_GoBack]ways init before using. That means init\n# before starting threads that use it.\" loses the mid-word bookmark\n# split and becomes one clean run: \"Initialize - always init...\".\n$initPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -match \"Initialize\" -and $t -match \"before using\") {\n        $initPara = $p\n        break\n    }\n}\n\nif ($initPara -ne $null) {\n    $fixedParaXml = '<w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"39\"/></w:numPr></w:pPr>' +\n        '<w:r><w:t>Initialize</w:t></w:r>' +\n        '<w:r><w:t xml:space=\"preserve\"> ' + $enDash + ' always </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> before using. That means </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>init</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> before starting threads that use it. </w:t></w:r>' +\n        '</w:p>'\n    $initPara.Range.InsertXML((New-FlatOpc $fixedParaXml))\n}\n\n# --- Change 2 ------------------------------------------------------------\n# Add a new indented paragraph right after the WICED_WAIT_FOREVER line\n# explaining that the \"forever\" timeout is really ~50 days; the\n# _GoBack bookmark (freed up by change 1) now wraps the point right\n# before \"days)\".\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"use WICED_WAIT_FOREVER if you don\" + $rsquo + \"t want a timeout\"\n$find.Forward = $true\n$find.Wrap = 0\n$found = $find.Execute()\n\nif ($found) {\n    $endPos = $find.Parent.End\n    $insertPoint = $d.Range($endPos, $endPos)\n    $newParaXml = '<w:p><w:pPr><w:ind w:left=\"720\"/></w:pPr>' +\n        '<w:r><w:tab/><w:t xml:space=\"preserve\">(it' + $rsquo + 's really 2^32 </w:t></w:r>' +\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>ms</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t xml:space=\"preserve\"> which is just under 50 </w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n        '<w:r><w:t>days)</w:t></w:r>' +\n        '</w:p>'\n    $insertPoint.InsertXML((New-FlatOpc $newParaXml))\n}\n"}
